# Add testNG reports feature (#10)
# Fix header typos and add a new "alerttext" column used to surface the
# "Customer added successfully" message raised by the testNG report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - correct the existing labels (typo fixes from the report work)
$ws.Range("A1").Value = "firstaame"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

# New column D: header + the alert text raised after adding a customer
$ws.Range("D1").Value = "alerttext"
$ws.Range("D2").Value = "Customer added successfully"

# Leave the current selection where the user last left off
$ws.Range("E12").Select() | Out-Null
